$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 45736.01041666666, 636),
    @(3, 45736.02083333334, 637),
    @(4, 45736.03125, 643),
    @(5, 45736.04166666666, 651),
    @(6, 45736.05208333334, 694),
    @(7, 45736.0625, 694),
    @(8, 45736.07291666666, 713),
    @(9, 45736.08333333334, 706),
    @(10, 45736.09375, 740),
    @(11, 45736.10416666666, 735),
    @(12, 45736.11458333334, 735),
    @(13, 45736.125, 731),
    @(14, 45736.13541666666, 723),
    @(15, 45736.14583333334, 718),
    @(16, 45736.15625, 716),
    @(17, 45736.16666666666, 719),
    @(18, 45736.17708333334, 692),
    @(19, 45736.1875, 703),
    @(20, 45736.19791666666, 675),
    @(21, 45736.20833333334, 679),
    @(22, 45736.21875, 721),
    @(23, 45736.22916666666, 696),
    @(24, 45736.23958333334, 688),
    @(25, 45736.25, 672),
    @(26, 45736.26041666666, 599),
    @(27, 45736.27083333334, 590),
    @(28, 45736.28125, 588),
    @(29, 45736.29166666666, 587),
    @(30, 45736.30208333334, 474),
    @(31, 45736.3125, 478),
    @(32, 45736.32291666666, 483),
    @(33, 45736.33333333334, 477),
    @(34, 45736.34375, 368),
    @(35, 45736.35416666666, 368),
    @(36, 45736.36458333334, 366),
    @(37, 45736.375, 367),
    @(38, 45736.38541666666, 364),
    @(39, 45736.39583333334, 364),
    @(40, 45736.40625, 365),
    @(41, 45736.41666666666, 365),
    @(42, 45736.42708333334, 346),
    @(43, 45736.4375, 346),
    @(44, 45736.44791666666, 345),
    @(45, 45736.45833333334, 344),
    @(46, 45736.46875, 334),
    @(47, 45736.47916666666, 333),
    @(48, 45736.48958333334, 332),
    @(49, 45736.5, 331),
    @(50, 45736.51041666666, 326),
    @(51, 45736.52083333334, 326),
    @(52, 45736.53125, 325),
    @(53, 45736.54166666666, 324),
    @(54, 45736.55208333334, 311),
    @(55, 45736.5625, 310),
    @(56, 45736.57291666666, 308),
    @(57, 45736.58333333334, 307),
    @(58, 45736.59375, 270),
    @(59, 45736.60416666666, 268),
    @(60, 45736.61458333334, 266),
    @(61, 45736.625, 264),
    @(62, 45736.63541666666, 214),
    @(63, 45736.64583333334, 210),
    @(64, 45736.65625, 207),
    @(65, 45736.66666666666, 204),
    @(66, 45736.67708333334, 136),
    @(67, 45736.6875, 133),
    @(68, 45736.69791666666, 131),
    @(69, 45736.70833333334, 129),
    @(70, 45736.71875, 95),
    @(71, 45736.72916666666, 94),
    @(72, 45736.73958333334, 94),
    @(73, 45736.75, 93),
    @(74, 45736.76041666666, 74),
    @(75, 45736.77083333334, 74),
    @(76, 45736.78125, 73),
    @(77, 45736.79166666666, 73),
    @(78, 45736.80208333334, 75),
    @(79, 45736.8125, 75),
    @(80, 45736.82291666666, 76),
    @(81, 45736.83333333334, 76),
    @(82, 45736.84375, 83),
    @(83, 45736.85416666666, 83),
    @(84, 45736.86458333334, 83),
    @(85, 45736.875, 83),
    @(86, 45736.88541666666, 87),
    @(87, 45736.89583333334, 89),
    @(88, 45736.90625, 89),
    @(89, 45736.91666666666, 90),
    @(90, 45736.92708333334, 97),
    @(91, 45736.9375, 97),
    @(92, 45736.94791666666, 97),
    @(93, 45736.95833333334, 98),
    @(94, 45736.96875, $null),
    @(95, 45736.97916666666, $null),
    @(96, 45736.98958333334, $null),
    @(97, 45737, $null)
)

foreach ($row in $data) {
    $r = $row[0]
    $dateVal = $row[1]
    $prodVal = $row[2]
    $ws.Cells.Item($r, 1).Value = $dateVal
    if ($null -ne $prodVal) {
        $ws.Cells.Item($r, 2).Value = $prodVal
    }
}
